$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2-31 down to 3-32)
$ws.Rows(2).Insert()

# Populate the newly inserted row with the Kraków entry
$ws.Range("A2").Value = "Kraków"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 50.062537797834501
$ws.Range("D2").Value = 19.937306291842098

# Update the active selection to match the author's final state
$null = $ws.Range("A3").Select()
